{"js": "// The commit removes the whole \"Source Code\" paragraph that listed the\n// working-directory files (## [1] \"Notofit.csv\" ... ## [9] \"Phase2_...csv\"),\n// right before the \"Supplementary methods for ...\" heading. The preceding\n// paragraph (ending in \"## This is DHARMa 0.4.7. ...\") is left untouched.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the paragraph that contains the distinctive \"Notofit.csv\" listing\n// text rather than relying on a fixed index (Word's own TOC field adds\n// extra paragraph-like items that can shift plain indices around).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"Notofit.csv\") !== -1 && text.indexOf(\"Phase2_dispersal_Notonecta.undulata.csv\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# The commit removes the whole \"Source Code\" paragraph that listed the\n# working-directory files (## [1] \"Notofit.csv\" ... ## [9] \"Phase2_...csv\"),\n# right before the \"Supplementary methods for ...\" heading. The preceding\n# paragraph (ending in \"## This is DHARMa 0.4.7. ...\") is left untouched.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -like \"*Notofit.csv*\" -and $t -like \"*Phase2_dispersal_Notonecta.undulata.csv*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
